$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coinranking price/volume refresh (GitHub Actions scheduled update).
# Numeric-looking Price values must stay text (matches original inline-string
# typing, e.g. "0.998" rather than the number 0.998), so we force the Text
# number format on those cells before writing the new value.

$ws.Range("D2").Value = '56.218.71'
$ws.Range("E2").Value = '  -0.96%  '
$ws.Range("D3").Value = '2.363.28'
$ws.Range("E3").Value = '  -1.32%  '
$ws.Range("E4").Value = '  +0.25%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '505.65'
$ws.Range("E5").Value = '  -0.37%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '132.63'
$ws.Range("E6").Value = '  -1.13%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.09%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.541'
$ws.Range("E8").Value = '  -2.22%  '
$ws.Range("D9").Value = '2.385.63'
$ws.Range("E9").Value = '  -0.73%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0964'
$ws.Range("E10").Value = '  -1.49%  '
$ws.Range("E11").Value = '  -0.55%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '4.83'
$ws.Range("E12").Value = '  +3.35%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.323'
$ws.Range("E13").Value = '  -4.25%  '
$ws.Range("D14").Value = '2.795.95'
$ws.Range("E14").Value = '  -0.73%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.85'
$ws.Range("E15").Value = '  +0.20%  '
$ws.Range("D16").Value = '56.199.69'
$ws.Range("E16").Value = '  -0.90%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.0000132'
$ws.Range("E17").Value = '  -0.72%  '
$ws.Range("D18").Value = '2.335.82'
$ws.Range("E18").Value = '  -2.97%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.00'
$ws.Range("E19").Value = '  -1.71%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.06'
$ws.Range("E20").Value = '  +0.23%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '310.82'
$ws.Range("E21").Value = '  +0.04%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.27'
$ws.Range("E22").Value = '  +0.21%  '
$ws.Range("E23").Value = '  -0.44%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.45'
$ws.Range("E24").Value = '  -0.24%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.996'
$ws.Range("E25").Value = '  -0.21%  '
$ws.Range("E26").Value = '  +0.33%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.148'
$ws.Range("E27").Value = '  -1.36%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.29'
$ws.Range("E28").Value = '  -0.74%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '170.69'
$ws.Range("E29").Value = '  -1.63%  '
$ws.Range("D30").Value = '0.0₃0717'
$ws.Range("E30").Value = '  -1.79%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.64'
$ws.Range("E31").Value = '  -1.31%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.84'
$ws.Range("E32").Value = '  -0.31%  '
$ws.Range("B33").Value = 'USDe'
$ws.Range("C33").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.998'
$ws.Range("E33").Value = '  -0.09%  '
$ws.Range("B34").Value = 'Fetch.AI'
$ws.Range("C34").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.09'
$ws.Range("E34").Value = '  -2.56%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.996'
$ws.Range("E35").Value = '  +0.21%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '17.77'
$ws.Range("E36").Value = '  -0.79%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.19'
$ws.Range("E37").Value = '  -0.95%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.877'
$ws.Range("E38").Value = '  +7.44%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.75'
$ws.Range("E39").Value = '  -2.60%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '36.52'
$ws.Range("E40").Value = '  -0.34%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.42'
$ws.Range("E41").Value = '  -1.38%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.378'
$ws.Range("E42").Value = '  +1.09%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.39'
$ws.Range("E43").Value = '  -0.33%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.98'
$ws.Range("E44").Value = '  +1.99%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '126.37'
$ws.Range("E45").Value = '  -4.84%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.562'
$ws.Range("E46").Value = '  -1.36%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0899'
$ws.Range("E47").Value = '  -1.16%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '244.76'
$ws.Range("E48").Value = '  -1.33%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0485'
$ws.Range("E49").Value = '  -0.37%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '17.08'
$ws.Range("E50").Value = '  +0.76%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0209'
$ws.Range("E51").Value = '  -0.93%  '
